$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144.41667
$ws.Range("J9").Value = 175.9
$ws.Range("L9").Value = 175.9
$ws.Range("N9").Value = -513.9

$ws.Range("H87").Value = 11004.238
$ws.Range("J87").Value = 11004.238
$ws.Range("L87").Value = 11004.238
$ws.Range("N87").Value = -13500.238

$ws.Range("H90").Value = 11004.238
$ws.Range("J90").Value = 11004.238
$ws.Range("L90").Value = 33012.714
$ws.Range("N90").Value = -45492.714

$ws.Range("H100").Value = 8800074
$ws.Range("I100").Value = 11905910
$ws.Range("J100").Value = 103732.3
$ws.Range("K100").Value = 11905910
$ws.Range("L100").Value = 103732.3
$ws.Range("M100").Value = -11905369
$ws.Range("N100").Value = -104814.3

$ws.Range("H113").Value = 11193.0625
$ws.Range("I113").Value = 4021.8
$ws.Range("J113").Value = 14452.728
$ws.Range("K113").Value = 4021.8
$ws.Range("L113").Value = 14452.728
$ws.Range("M113").Value = -767.8000000000002
$ws.Range("N113").Value = -20960.728

$ws.Range("H129").Value = 1214.6786
$ws.Range("I129").Value = 633.5
$ws.Range("J129").Value = 1311.5416
$ws.Range("K129").Value = 1900.5
$ws.Range("L129").Value = 3934.6248
$ws.Range("M129").Value = 3099.5
$ws.Range("N129").Value = -13934.6248

$ws.Range("H138").Value = 1832.258
$ws.Range("I138").Value = 1135.4478
$ws.Range("J138").Value = 3627.8845
$ws.Range("K138").Value = 3406.3434
$ws.Range("L138").Value = 10883.6535
$ws.Range("M138").Value = 1733.6566
$ws.Range("N138").Value = -21163.6535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 768.68085
$ws.Range("I61").Value = 565.1111
$ws.Range("J61").Value = 1434.909
$ws.Range("K61").Value = 565.1111
$ws.Range("L61").Value = 1434.909
$ws.Range("M61").Value = -353.1111
$ws.Range("N61").Value = -1858.909

$ws.Range("H132").Value = 2347.3157
$ws.Range("I132").Value = 1128.0555
$ws.Range("J132").Value = 3444.65
$ws.Range("K132").Value = 3384.1665
$ws.Range("L132").Value = 10333.95
$ws.Range("M132").Value = -854.1664999999998
$ws.Range("N132").Value = -15393.95

$ws.Range("H136").Value = 768.68085
$ws.Range("I136").Value = 565.1111
$ws.Range("J136").Value = 1434.909
$ws.Range("K136").Value = 1695.3333
$ws.Range("L136").Value = 4304.727000000001
$ws.Range("M136").Value = 854.6667000000002
$ws.Range("N136").Value = -9404.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7606.278
$ws.Range("I20").Value = 1866.25
$ws.Range("J20").Value = 19086.334
$ws.Range("K20").Value = 1866.25
$ws.Range("L20").Value = 19086.334
$ws.Range("M20").Value = -1619.25
$ws.Range("N20").Value = -19580.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2533.5
$ws.Range("I31").Value = 1563.1613
$ws.Range("J31").Value = 3503.8386
$ws.Range("K31").Value = 1563.1613
$ws.Range("L31").Value = 3503.8386
$ws.Range("M31").Value = -1268.1613
$ws.Range("N31").Value = -4093.8386

$ws.Range("H34").Value = 2533.5
$ws.Range("I34").Value = 1563.1613
$ws.Range("J34").Value = 3503.8386
$ws.Range("K34").Value = 1563.1613
$ws.Range("L34").Value = 3503.8386
$ws.Range("M34").Value = -1361.1613
$ws.Range("N34").Value = -3907.8386

$ws.Range("H122").Value = 127842.25
$ws.Range("I122").Value = 167956.33
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 503868.99
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -501418.99
$ws.Range("N122").Value = -27400

$ws.Range("H132").Value = 1536.234
$ws.Range("I132").Value = 961.9211
$ws.Range("J132").Value = 3961.111
$ws.Range("K132").Value = 2885.7633
$ws.Range("L132").Value = 11883.333
$ws.Range("M132").Value = -355.7633000000001
$ws.Range("N132").Value = -16943.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2831.3914
$ws.Range("I109").Value = 542.3333
$ws.Range("J109").Value = 3174.75
$ws.Range("K109").Value = 1626.9999
$ws.Range("L109").Value = 9524.25
$ws.Range("M109").Value = -586.9999
$ws.Range("N109").Value = -11604.25

$ws.Range("H110").Value = 3333.3333
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 4500
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 13500
$ws.Range("M110").Value = 1090
$ws.Range("N110").Value = -21680

$ws.Range("H112").Value = 1368.5294
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 1384.2424
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 4152.7272
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -6368.7272

$ws.Range("H114").Value = 1060
$ws.Range("J114").Value = 2266.6667
$ws.Range("L114").Value = 6800.000100000001
$ws.Range("N114").Value = -13308.0001

$ws.Range("H115").Value = 2776.2222
$ws.Range("I115").Value = 659.3333
$ws.Range("K115").Value = 1977.9999
$ws.Range("M115").Value = -802.9999

$ws.Range("H121").Value = 32139.188
$ws.Range("I121").Value = 200278.6
$ws.Range("J121").Value = 1002.2593
$ws.Range("K121").Value = 600835.8
$ws.Range("L121").Value = 3006.7779
$ws.Range("M121").Value = -599525.8
$ws.Range("N121").Value = -5626.7779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4998.08
$ws.Range("I70").Value = 4385.2646
$ws.Range("J70").Value = 6300.3125
$ws.Range("K70").Value = 4385.2646
$ws.Range("L70").Value = 6300.3125
$ws.Range("M70").Value = -4115.2646
$ws.Range("N70").Value = -6840.3125

$ws.Range("H73").Value = 4998.08
$ws.Range("I73").Value = 4385.2646
$ws.Range("J73").Value = 6300.3125
$ws.Range("K73").Value = 4385.2646
$ws.Range("L73").Value = 6300.3125
$ws.Range("M73").Value = -3449.2646
$ws.Range("N73").Value = -8172.3125

$ws.Range("H102").Value = 3116.6667
$ws.Range("I102").Value = 2675
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 2675
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -1053
$ws.Range("N102").Value = -7244

$ws.Range("H126").Value = 1838.1666
$ws.Range("I126").Value = 1725.8
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5177.4
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -2707.4
$ws.Range("N126").Value = -12140

$ws.Range("H132").Value = 1583.7059
$ws.Range("I132").Value = 1284.55
$ws.Range("J132").Value = 2671.5454
$ws.Range("K132").Value = 3853.65
$ws.Range("L132").Value = 8014.6362
$ws.Range("M132").Value = -1323.65
$ws.Range("N132").Value = -13074.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1876.6666
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1876.6666
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 1876.6666
$ws.Range("N68").Value = -3374.6666
$ws.Range("M68").Value = $null

$ws.Range("H71").Value = 1876.6666
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1876.6666
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 9383.333000000001
$ws.Range("N71").Value = -16871.333
$ws.Range("M71").Value = $null

$ws.Range("H132").Value = 6025.0435
$ws.Range("I132").Value = 6186.9316
$ws.Range("J132").Value = 5740.12
$ws.Range("K132").Value = 18560.7948
$ws.Range("L132").Value = 17220.36
$ws.Range("M132").Value = -16030.7948
$ws.Range("N132").Value = -22280.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4382.091
$ws.Range("I62").Value = 3401
$ws.Range("J62").Value = 4750
$ws.Range("K62").Value = 3401
$ws.Range("L62").Value = 4750
$ws.Range("M62").Value = -2777
$ws.Range("N62").Value = -5998

$ws.Range("H65").Value = 4382.091
$ws.Range("I65").Value = 3401
$ws.Range("J65").Value = 4750
$ws.Range("K65").Value = 17005
$ws.Range("L65").Value = 23750
$ws.Range("M65").Value = -13885
$ws.Range("N65").Value = -29990

$ws.Range("H132").Value = 1324.4166
$ws.Range("I132").Value = 1234.74
$ws.Range("J132").Value = 1772.8
$ws.Range("K132").Value = 3704.22
$ws.Range("L132").Value = 5318.4
$ws.Range("M132").Value = -1174.22
$ws.Range("N132").Value = -10378.4

$ws.Range("H136").Value = 1377.4259
$ws.Range("I136").Value = 517.8333
$ws.Range("J136").Value = 8254.166999999999
$ws.Range("K136").Value = 1553.4999
$ws.Range("L136").Value = 24762.501
$ws.Range("M136").Value = 996.5001
$ws.Range("N136").Value = -29862.501
